$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.585.28"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.158.99"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.35"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.93"
$ws.Range("E6").Value = "  -4.48%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.155.92"
$ws.Range("E8").Value = "  -4.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("E9").Value = "  -2.43%  "
$ws.Range("E10").Value = "  -5.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -3.71%  "
$ws.Range("E13").Value = "  -5.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.74"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.680.76"
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.160.71"
$ws.Range("E17").Value = "  -4.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.543.08"
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.56"
$ws.Range("E19").Value = "  -4.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "453.16"
$ws.Range("E20").Value = "  -5.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.97"
$ws.Range("E21").Value = "  -1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("E23").Value = "  -5.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.34"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.54"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.71"
$ws.Range("E29").Value = "  -5.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.73"
$ws.Range("E30").Value = "  -6.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.02"
$ws.Range("E31").Value = "  -6.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.19"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.38"
$ws.Range("E34").Value = "  -6.87%  "
$ws.Range("E35").Value = "  -6.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.82"
$ws.Range("E36").Value = "  -3.09%  "
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0703"
$ws.Range("E38").Value = "  -5.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0386"
$ws.Range("E39").Value = "  -3.64%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "402.08"
$ws.Range("E40").Value = "  -7.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.03"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.111"
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.786.50"
$ws.Range("E44").Value = "  -9.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.249"
$ws.Range("E45").Value = "  -6.16%  "

# Row 46: was Fetch.AI -> now USDe
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47: was USDe -> now Fetch.AI
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("E47").Value = "  -3.68%  "

# Row 48: Monero, D/E only
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.29"
$ws.Range("E48").Value = "  -2.52%  "

# Row 49: was InjectiveProtocol -> now Arweave
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.80"
$ws.Range("E49").Value = "  -3.45%  "

# Row 50: was Arweave -> now InjectiveProtocol
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.27"
$ws.Range("E50").Value = "  -4.07%  "

# Row 51: Stellar, E only
$ws.Range("E51").Value = "  -3.31%  "
